# Update betting-odds values on Sheet1 to match the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 ---
$ws.Range("G5").Value  = 1.9
$ws.Range("H5").Value  = 3.6
$ws.Range("I5").Value  = 3.8
$ws.Range("J5").Value  = 2.63
$ws.Range("L5").Value  = 4.75
$ws.Range("Q5").Value  = 2.3
$ws.Range("R5").Value  = 1.6
$ws.Range("AG5").Value = 8.5
$ws.Range("AM5").Value = 800
$ws.Range("AO5").Value = 11
$ws.Range("AQ5").Value = 41
$ws.Range("AX5").Value = 23

# --- Row 10 ---
$ws.Range("G10").Value  = 3.8
$ws.Range("H10").Value  = 3.3
$ws.Range("I10").Value  = 1.9
$ws.Range("J10").Value  = 4.75
$ws.Range("L10").Value  = 2.63
$ws.Range("S10").Value  = 1.44
$ws.Range("T10").Value  = 2.63
$ws.Range("Y10").Value  = 15
$ws.Range("AN10").Value = 6
$ws.Range("AT10").Value = 2.63
$ws.Range("AW10").Value = 3.75

# --- Row 12 ---
$ws.Range("G12").Value  = 2.3
$ws.Range("I12").Value  = 2.88
$ws.Range("W12").Value  = 7.5
$ws.Range("AN12").Value = 4.33
$ws.Range("AW12").Value = 5

# --- Row 15 ---
$ws.Range("W15").Value  = 6.5
$ws.Range("AD15").Value = 6
$ws.Range("AF15").Value = 67
$ws.Range("AU15").Value = 9
$ws.Range("AW15").Value = 5

# --- Row 19 ---
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = 1.85

# --- Row 24 ---
$ws.Range("H24").Value  = 3.75
$ws.Range("I24").Value  = 4.55
$ws.Range("J24").Value  = 2.18
$ws.Range("L24").Value  = 4.85
$ws.Range("M24").Value  = 1.05
$ws.Range("N24").Value  = 8
$ws.Range("P24").Value  = 3.65
$ws.Range("T24").Value  = 2.9
$ws.Range("X24").Value  = 8.25
$ws.Range("AB24").Value = 23
$ws.Range("AC24").Value = 8
$ws.Range("AD24").Value = 7.5
$ws.Range("AH24").Value = 28
$ws.Range("AJ24").Value = 80
$ws.Range("AK24").Value = 45
$ws.Range("AT24").Value = 2.9
$ws.Range("AW24").Value = 6.4
$ws.Range("AX24").Value = 26
$ws.Range("BB24").Value = 400
